$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Week 1 assignment link (J2): install-r-github -> quarto-notes
$ws.Range("J2").Value = "/assignment/01-assignment-quarto-notes"

# Rearrange the weekly "title" column (F) per the new schedule.
# Order matches the order new strings were introduced so the shared-string
# table comes out in the same sequence as the target workbook.
$ws.Range("F3").Value  = "No class (Labor Day)"
$ws.Range("F5").Value  = "Using R to look at data"
$ws.Range("F11").Value = "Working with models"
$ws.Range("F13").Value = "Functional programming patterns"
$ws.Range("F7").Value  = "Ingesting and cleaning data"
$ws.Range("F6").Value  = "Tidy data and dplyr"
$ws.Range("F14").Value = "Build systems, environments, and packages"
$ws.Range("F8").Value  = "Better tables, better graphs"

$ws.Range("F2").Value  = "Big Picture: Doing your work properly"
$ws.Range("F4").Value  = "The file system; the shell; the terminal"
$ws.Range("F9").Value  = "No class (Fall break)"
$ws.Range("F10").Value = "Version Control: git and GitHub"
$ws.Range("F12").Value = "Databases and APIs"
# F15 text is unchanged ("Leveraging Minions: ...") - leave it untouched so it
# keeps reusing its existing shared-string entry instead of being rewritten.

# Move the selection cursor to F6 (matches the saved view state in the diff)
$ws.Range("F6").Select()
